# Updated cryptos list - applies new Price (D) and Volume(1h) (E) values
# per-row, matching the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "26.232.23"
$ws.Cells.Item(2, 5).Value = "  +0.55%  "
$ws.Cells.Item(3, 4).Value = "1.659.00"
$ws.Cells.Item(3, 5).Value = "  +0.13%  "
$ws.Cells.Item(4, 5).Value = "  +0.65%  "
$ws.Cells.Item(5, 4).Value = "'218.23"
$ws.Cells.Item(5, 5).Value = "  -0.11%  "
$ws.Cells.Item(6, 5).Value = "  +0.11%  "
$ws.Cells.Item(7, 4).Value = "'1.008"
$ws.Cells.Item(8, 4).Value = "'0.2635"
$ws.Cells.Item(8, 5).Value = "  +0.70%  "
$ws.Cells.Item(9, 5).Value = "  +0.17%  "
$ws.Cells.Item(10, 4).Value = "'20.45"
$ws.Cells.Item(10, 5).Value = "  +0.08%  "
$ws.Cells.Item(11, 4).Value = "'0.07828"
$ws.Cells.Item(11, 5).Value = "  +0.81%  "
$ws.Cells.Item(12, 4).Value = "'4.541"
$ws.Cells.Item(12, 5).Value = "  +1.04%  "
$ws.Cells.Item(13, 4).Value = "1.671.50"
$ws.Cells.Item(13, 5).Value = "  +1.37%  "
$ws.Cells.Item(14, 4).Value = "1.887.24"
$ws.Cells.Item(14, 5).Value = "  +0.18%  "
$ws.Cells.Item(15, 4).Value = "'0.5522"
$ws.Cells.Item(15, 5).Value = "  +0.72%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8174"
$ws.Cells.Item(16, 5).Value = "  +0.34%  "
$ws.Cells.Item(17, 4).Value = "'65.53"
$ws.Cells.Item(17, 5).Value = "  +0.50%  "
$ws.Cells.Item(18, 4).Value = "'1.009"
$ws.Cells.Item(18, 5).Value = "  +0.68%  "
$ws.Cells.Item(19, 4).Value = "'4.639"
$ws.Cells.Item(19, 5).Value = "  +1.97%  "
$ws.Cells.Item(20, 4).Value = "'192.07"
$ws.Cells.Item(20, 5).Value = "  -0.58%  "
$ws.Cells.Item(21, 5).Value = "  +1.07%  "
$ws.Cells.Item(22, 4).Value = "'6.043"
$ws.Cells.Item(22, 5).Value = "  +0.35%  "
$ws.Cells.Item(23, 4).Value = "'1.010"
$ws.Cells.Item(23, 5).Value = "  +0.61%  "
$ws.Cells.Item(24, 4).Value = "'144.20"
$ws.Cells.Item(24, 5).Value = "  +2.92%  "
$ws.Cells.Item(25, 4).Value = "'0.1221"
$ws.Cells.Item(25, 5).Value = "  -1.81%  "
$ws.Cells.Item(26, 4).Value = "'7.215"
$ws.Cells.Item(26, 5).Value = "  -0.78%  "
$ws.Cells.Item(27, 4).Value = "'16.10"
$ws.Cells.Item(27, 5).Value = "  -0.41%  "
$ws.Cells.Item(29, 4).Value = "'0.05861"
$ws.Cells.Item(29, 5).Value = "  -1.42%  "
$ws.Cells.Item(30, 4).Value = "'1.277"
$ws.Cells.Item(30, 5).Value = "  +0.05%  "
$ws.Cells.Item(31, 4).Value = "'3.575"
$ws.Cells.Item(31, 5).Value = "  +1.88%  "
$ws.Cells.Item(32, 5).Value = "  +1.62%  "
$ws.Cells.Item(33, 4).Value = "'1.612"
$ws.Cells.Item(33, 5).Value = "  +3.36%  "
$ws.Cells.Item(34, 4).Value = "'0.9596"
$ws.Cells.Item(34, 5).Value = "  +1.13%  "
$ws.Cells.Item(35, 4).Value = "'2.818"
$ws.Cells.Item(35, 5).Value = "  +1.76%  "
$ws.Cells.Item(36, 4).Value = "'2.421"
$ws.Cells.Item(36, 5).Value = "  +0.41%  "
$ws.Cells.Item(37, 4).Value = "'0.5803"
$ws.Cells.Item(37, 5).Value = "  +2.84%  "
$ws.Cells.Item(38, 4).Value = "'0.01603"
$ws.Cells.Item(38, 5).Value = "  -0.38%  "
$ws.Cells.Item(39, 4).Value = "'5.900"
$ws.Cells.Item(39, 5).Value = "  +1.08%  "
$ws.Cells.Item(40, 4).Value = "'0.8519"
$ws.Cells.Item(40, 5).Value = "  +0.92%  "
$ws.Cells.Item(41, 5).Value = "  +0.58%  "
$ws.Cells.Item(42, 4).Value = "1.044.91"
$ws.Cells.Item(42, 5).Value = "  +3.24%  "
$ws.Cells.Item(43, 4).Value = "'103.87"
$ws.Cells.Item(43, 5).Value = "  +2.37%  "
$ws.Cells.Item(44, 4).Value = "1.800.33"
$ws.Cells.Item(44, 5).Value = "  +0.05%  "
$ws.Cells.Item(45, 4).Value = "'57.06"
$ws.Cells.Item(45, 5).Value = "  +0.02%  "
$ws.Cells.Item(46, 5).Value = "  +3.19%  "
$ws.Cells.Item(47, 5).Value = "  +1.17%  "
$ws.Cells.Item(48, 4).Value = "'0.4369"
$ws.Cells.Item(48, 5).Value = "  +1.95%  "
$ws.Cells.Item(49, 4).Value = "'7.934"
$ws.Cells.Item(49, 5).Value = "  +2.85%  "
$ws.Cells.Item(50, 4).Value = "'0.05162"
$ws.Cells.Item(50, 5).Value = "  +0.12%  "
$ws.Cells.Item(51, 4).Value = "'1.432"
$ws.Cells.Item(51, 5).Value = "  -2.50%  "
